# 7/10/20: ran scrape and analysis
# Adds a new "BF" column (date 7/10/2020, serial 44022) to each of the four
# data sheets, mirroring the existing BE column's per-row structure, and
# restores each sheet's on-screen selection to where the author left it.

$wb = $excel.ActiveWorkbook

$data = @{
    "Facilities" = @{ 4 = 17;  5 = 114; 6 = 12; 7 = 118; 8 = 35;  9 = 29; 10 = 12; 11 = 7; 12 = 344 }
    "Cases"      = @{ 4 = 40;  5 = 364; 6 = 18; 7 = 322; 8 = 139; 9 = 13; 10 = 59; 11 = 7; 12 = 962 }
    "Fatalities" = @{ 4 = 2;   5 = 83;  6 = 5;  7 = 39;  8 = 20;  9 = 5;  10 = 1;  11 = 2; 12 = 157 }
    "Recoveries" = @{ 4 = 19;  5 = 133; 6 = 11; 7 = 114; 8 = 74;  9 = 6;  10 = 2;  11 = 1; 12 = 360 }
}

# Final on-sheet selection per tab, matching where the author left the
# cursor after pasting in the new day's numbers.
$selections = @{
    "Cases"      = "BE17"
    "Fatalities" = "BF4:BF12"
    "Recoveries" = "BF4:BF12"
    "Facilities" = "A2"
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New date header in BF3, continuing the daily series from BE3 (44021).
    $ws.Range("BF3").Value = 44022

    $rows = $data[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 58).Value = $rows[$r]
    }
}

# Apply selections; do "Facilities" last so it remains the active tab,
# matching the saved workbook's original tabSelected state.
foreach ($sheetName in @("Cases", "Fatalities", "Recoveries", "Facilities")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($selections[$sheetName]).Select()
}
